$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift the existing data rows (2-6) down by one to make room for the new
# "KL - Richmond" row, dropping the last row ("Default Store Name") off the
# bottom (the sheet stays A1:B7).
for ($r = 6; $r -ge 2; $r--) {
    $srcA = $ws.Cells.Item($r, 1).Value2
    $srcB = $ws.Cells.Item($r, 2).Value2
    $ws.Cells.Item($r + 1, 1).Value = $srcA
    $ws.Cells.Item($r + 1, 2).Value = $srcB
}

# Write the new first data row
$ws.Cells.Item(2, 1).Value = "KL - Richmond"
$ws.Cells.Item(2, 2).Value = "KL - Richmond"
